# Add a new "Player Info" worksheet before "ODI Batting", populate it,
# and rewrite the MATCH_CARD_LINK column on "ODI Batting" into MATCH_CODE.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "Player Info" sheet, placed before "ODI Batting" --
$originalBatting = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($originalBatting)
$playerInfo.Name = "Player Info"

# Re-resolve "ODI Batting" by name now that the sheet collection changed -
# a reference captured before the Add() call can silently track the old
# positional index instead of staying bound to the sheet itself.
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# --- 2. Populate "Player Info" -------------------------------------------
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Reuse the exact header formatting already used on "ODI Batting" (bold,
# centered, bordered) instead of re-deriving it property by property.
$battingSheet.Range("A1:D1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$playerInfo.Range("A2").Value = "6037"
$playerInfo.Range("B2").Value = "John Andrew Simpson"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Does Not Bowl | Unknown"

# --- 3. Update "ODI Batting": rename header & shorten link to match code -
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingSheet.Range("D2").Value = "4472"
$battingSheet.Range("D3").Value = "4473"
$battingSheet.Range("D4").Value = "4476"
